$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the "Total cost" summary values in row 20 (Less-or-equal-to-5 H-index
#        group) so every count collapses to 0 (keeping the reported upper bound). ---
$ws.Range("B20").Value = "0 (0, 795)"
$ws.Range("C20").Value = "0 (0, 139)"
$ws.Range("D20").Value = "0 (0, 88)"
$ws.Range("E20").Value = "1,197 (204, 3,376)"
$ws.Range("F20").Value = "3,024 (1,548, 5,449)"
$ws.Range("G20").Value = "0 (0, 1,901)"

# --- 2. Row 21 (the "Unknown" H-Index count row) is being removed entirely. Before
#        deleting it, copy its formatting onto row 20 so row 20 inherits the thin
#        bottom border that previously marked the end of the data block. ---
$ws.Range("A21:P21").Copy()
$ws.Range("A20:P20").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 3. Delete row 21 outright - rows below shift up (old 22/23/24 become 21/22/23),
#        and the merged-cell ranges for the footnote rows renumber automatically. ---
$ws.Rows.Item(21).Delete()
